$d = $word.ActiveDocument

# Merge the split runs in the Title paragraph ("Answers:" / " " / "Introduction" /
# " " / "to" / " " / "factorization") into a single run of text. wdReplaceOne (1)
# is used so only this (first/only) occurrence is touched.
$d.Content.Find.Execute("Answers: Introduction to factorization", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Answers: Introduction to factorization", 1)

# Merge the split runs in the Author paragraph ("Millie" / " " / "Pike") into a
# single run of text. wdReplaceOne (1) ensures the later, unrelated "Millie Pike"
# mention in the version-history paragraph is left untouched.
$d.Content.Find.Execute("Millie Pike", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Millie Pike", 1)

# Merge the split runs in the Abstract paragraph into a single run of text.
$d.Content.Find.Execute("Answers to questions relating to the guide on introduction to factorisation.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Answers to questions relating to the guide on introduction to factorisation.", 1)
